# Update the "Clases" (class intervals) labels in the frequency table
# from parenthesis notation "(a, b)" to bracket notation "[a, b]" and
# correct a few rounding values, per the commit:
# "Se termina la practica Tabla de Frecuencias ..."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabla de frecuencias")

$ws.Range("A4").Value = "[20.7, 23.51]"
$ws.Range("A5").Value = "[23.51, 26.33]"
$ws.Range("A6").Value = "[26.33, 29.14]"
$ws.Range("A7").Value = "[29.14, 31.96]"
$ws.Range("A8").Value = "[31.96, 34.77]"
$ws.Range("A9").Value = "[34.77, 37.59]"
$ws.Range("A10").Value = "[37.59, 40.4]"
